$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates: volume/number and week-covering dates ---
$ws.Range("A8").Value = "Volume 30   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/24/2023  Through  7/30/2023"

# --- Weekly crime-complaints table (rows 14-30): refreshed figures ---
$ws.Range("D14").Value = 1
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 10
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 100
$ws.Range("M14").Value = -37.5
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = -25
$ws.Range("I15").Value = 26
$ws.Range("J15").Value = 33
$ws.Range("K15").Value = -21.212121212121
$ws.Range("L15").Value = 13.043478260869
$ws.Range("M15").Value = 4
$ws.Range("N15").Value = -36.585365853658
$ws.Range("C16").Value = 16
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 60
$ws.Range("F16").Value = 53
$ws.Range("G16").Value = 53
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 271
$ws.Range("J16").Value = 256
$ws.Range("K16").Value = 5.859375
$ws.Range("L16").Value = 47.282608695652
$ws.Range("M16").Value = 11.522633744856
$ws.Range("N16").Value = -66.625615763546
$ws.Range("C17").Value = 17
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = -5.555555555555
$ws.Range("F17").Value = 64
$ws.Range("G17").Value = 72
$ws.Range("H17").Value = -11.111111111111
$ws.Range("I17").Value = 470
$ws.Range("J17").Value = 431
$ws.Range("K17").Value = 9.048723897911
$ws.Range("L17").Value = 14.355231143552
$ws.Range("M17").Value = 91.836734693877
$ws.Range("N17").Value = -9.266409266409
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 28
$ws.Range("H18").Value = -14.285714285714
$ws.Range("I18").Value = 170
$ws.Range("J18").Value = 174
$ws.Range("K18").Value = -2.298850574712
$ws.Range("L18").Value = 16.438356164383
$ws.Range("M18").Value = -13.705583756345
$ws.Range("N18").Value = -84.942426926483
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 31.25
$ws.Range("F19").Value = 65
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = 51.162790697674
$ws.Range("I19").Value = 434
$ws.Range("J19").Value = 426
$ws.Range("K19").Value = 1.8779342723
$ws.Range("L19").Value = 54.448398576512
$ws.Range("M19").Value = 158.333333333333
$ws.Range("N19").Value = 37.777777777777
$ws.Range("C20").Value = 18
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 125
$ws.Range("F20").Value = 48
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = 84.615384615384
$ws.Range("I20").Value = 372
$ws.Range("J20").Value = 270
$ws.Range("K20").Value = 37.777777777777
$ws.Range("L20").Value = 27.835051546391
$ws.Range("M20").Value = 96.825396825396
$ws.Range("N20").Value = -60.128617363344
$ws.Range("C21").Value = 81
$ws.Range("D21").Value = 64
$ws.Range("E21").Value = 26.5625
$ws.Range("F21").Value = 262
$ws.Range("G21").Value = 232
$ws.Range("H21").Value = 12.931034482758
$ws.Range("I21").Value = 1753
$ws.Range("J21").Value = 1600
$ws.Range("K21").Value = 9.5625
$ws.Range("L21").Value = 30.723340790454
$ws.Range("M21").Value = 61.865189289012
$ws.Range("N21").Value = -53.464295195115
$ws.Range("C22").Value = 2
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = 100
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 15
$ws.Range("J22").Value = 21
$ws.Range("K22").Value = -28.571428571428
$ws.Range("L22").Value = 66.666666666666
$ws.Range("M22").Value = -16.666666666666
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = 36.363636363636
$ws.Range("I23").Value = 62
$ws.Range("J23").Value = 71
$ws.Range("K23").Value = -12.676056338028
$ws.Range("L23").Value = -3.125
$ws.Range("M23").Value = 67.567567567567
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -25.806451612903
$ws.Range("F24").Value = 90
$ws.Range("G24").Value = 99
$ws.Range("H24").Value = -9.090909090909
$ws.Range("I24").Value = 788
$ws.Range("J24").Value = 875
$ws.Range("K24").Value = -9.942857142857
$ws.Range("L24").Value = 29.818780889621
$ws.Range("M24").Value = 95.049504950495
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -38.095238095238
$ws.Range("F25").Value = 80
$ws.Range("G25").Value = 100
$ws.Range("H25").Value = -20
$ws.Range("I25").Value = 560
$ws.Range("J25").Value = 551
$ws.Range("K25").Value = 1.633393829401
$ws.Range("L25").Value = 12.903225806451
$ws.Range("M25").Value = -4.273504273504
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 7
$ws.Range("G26").Value = 10
$ws.Range("H26").Value = -30
$ws.Range("I26").Value = 37
$ws.Range("J26").Value = 54
$ws.Range("K26").Value = -31.481481481481
$ws.Range("L26").Value = 2.777777777777
$ws.Range("C27").Value = 4
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 300
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 16.666666666666
$ws.Range("I27").Value = 46
$ws.Range("J27").Value = 37
$ws.Range("K27").Value = 24.324324324324
$ws.Range("L27").Value = 17.948717948717
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C28").NumberFormat = "General"
$ws.Range("D28").Value = 2
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -33.333333333333
$ws.Range("J28").Value = 32
$ws.Range("K28").Value = -15.625
$ws.Range("M28").Value = -35.714285714285
$ws.Range("N28").Value = -65.822784810126
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C29").NumberFormat = "General"
$ws.Range("D29").Value = 2
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = -20
$ws.Range("J29").Value = 30
$ws.Range("K29").Value = -23.333333333333
$ws.Range("M29").Value = -32.35294117647
$ws.Range("N29").Value = -68.918918918918